# Azure Bootcamp 2019.xlsx edit:
#  - rename Sheet1 -> "Project Fields"
#  - add two new sheets: "Project Allocation" and "Project Seasonality"
#  - populate both new sheets with their tables / formulas
#  - restore per-sheet selections / active sheet

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename the original sheet.
# ---------------------------------------------------------------------------
$wsFields = $wb.Worksheets.Item(1)
$wsFields.Name = "Project Fields"

# ---------------------------------------------------------------------------
# 2. Insert the two new sheets right after "Project Fields".
#    Inserting "Project Seasonality" first and "Project Allocation" second
#    (both anchored After $wsFields) reproduces the original authoring
#    order: Seasonality picks up the lower internal sheetId, Allocation the
#    higher one, while ending up positioned *before* Seasonality in the tab
#    strip - exactly matching the target tab order/sheetId/r:id layout.
# ---------------------------------------------------------------------------
$tmp = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsFields)
$tmp.Name = "Project Seasonality"
$tmp = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsFields)
$tmp.Name = "Project Allocation"

# Re-resolve stable handles by name - worksheet handles obtained before an
# Add() can silently repoint to a different sheet once indices shift.
$wsFields = $wb.Worksheets.Item("Project Fields")
$wsAllocation = $wb.Worksheets.Item("Project Allocation")
$wsSeasonality = $wb.Worksheets.Item("Project Seasonality")

# ---------------------------------------------------------------------------
# 3. "Project Seasonality" sheet header row + row labels: Project Name /
#    Mon 1..Mon 12. Written before the "Project Allocation" sheet's header
#    so the new shared-string entries land in the same order the original
#    authoring session produced them in (Mon 1..Mon 12, then Allocation %,
#    then the trailing "x" marker last) - see step 5 below.
# ---------------------------------------------------------------------------
$wsSeasonality.Range("A1").Value = "Project Name"
$wsSeasonality.Range("B1").Value = "Mon 1"
$wsSeasonality.Range("C1").Value = "Mon 2"
$wsSeasonality.Range("D1").Value = "Mon 3"
$wsSeasonality.Range("E1").Value = "Mon 4"
$wsSeasonality.Range("F1").Value = "Mon 5"
$wsSeasonality.Range("G1").Value = "Mon 6"
$wsSeasonality.Range("H1").Value = "Mon 7"
$wsSeasonality.Range("I1").Value = "Mon 8"
$wsSeasonality.Range("J1").Value = "Mon 9"
$wsSeasonality.Range("K1").Value = "Mon 10"
$wsSeasonality.Range("L1").Value = "Mon 11"
$wsSeasonality.Range("M1").Value = "Mon 12"

$wsSeasonality.Range("A2").Value = "Krypton"
$wsSeasonality.Range("A3").Value = "Clingon"
$wsSeasonality.Range("A4").Value = "Jupiter"
$wsSeasonality.Range("A5").Value = "Tesla"
$wsSeasonality.Range("A6").Value = "Venus"

# ---------------------------------------------------------------------------
# 4. "Project Allocation" sheet: Project Name / Allocation % table.
# ---------------------------------------------------------------------------
$wsAllocation.Range("A1").Value = "Project Name"
$wsAllocation.Range("B1").Value = "Allocation %"

$wsAllocation.Range("A2").Value = "Krypton"
$wsAllocation.Range("A3").Value = "Clingon"
$wsAllocation.Range("A4").Value = "Jupiter"
$wsAllocation.Range("A5").Value = "Tesla"
$wsAllocation.Range("A6").Value = "Venus"

$wsAllocation.Range("B2").Value = 0.25
$wsAllocation.Range("B3").Value = 0.5
$wsAllocation.Range("B4").Value = 0.15
$wsAllocation.Range("B5").Value = 0.05
$wsAllocation.Range("B6").Value = 0.05

# Number formats for the data block.
$wsSeasonality.Range("B2:M6").NumberFormat = "0.0000"
$wsSeasonality.Range("N2:N6").NumberFormat = '_(* #,##0.000_);_(* \(#,##0.000\);_(* "-"??_);_(@_)'
$wsSeasonality.Range("B7:M7").NumberFormat = "0.00"
$wsSeasonality.Range("B8").NumberFormat = "0.000"
$wsSeasonality.Range("C8:M8").NumberFormat = "0.00"

# Row 2 - Krypton: even 1/12 split (entered as literal decimals).
$wsSeasonality.Range("B2").Value = 0.083333333
$wsSeasonality.Range("C2").Value = 0.083333333
$wsSeasonality.Range("D2").Value = 0.083333333
$wsSeasonality.Range("E2").Value = 0.083333333
$wsSeasonality.Range("F2").Value = 0.083333333
$wsSeasonality.Range("G2").Value = 0.083333333
$wsSeasonality.Range("H2").Value = 0.083333333
$wsSeasonality.Range("I2").Value = 0.083333333
$wsSeasonality.Range("J2").Value = 0.083333333
$wsSeasonality.Range("K2").Value = 0.083333333
$wsSeasonality.Range("L2").Value = 0.083333333
$wsSeasonality.Range("M2").Value = 0.083333333
$wsSeasonality.Range("N2").Formula = "=SUM(B2:M2)"

# Row 3 - Clingon.
$wsSeasonality.Range("B3").Value = 0.03
$wsSeasonality.Range("C3").Value = 0.035
$wsSeasonality.Range("D3").Value = 0.035
$wsSeasonality.Range("E3").Value = 0.055
$wsSeasonality.Range("F3").Value = 0.055
$wsSeasonality.Range("G3").Value = 0.1
$wsSeasonality.Range("H3").Value = 0.07
$wsSeasonality.Range("I3").Value = 0.07
$wsSeasonality.Range("J3").Value = 0.1
$wsSeasonality.Range("K3").Value = 0.1
$wsSeasonality.Range("L3").Value = 0.15
$wsSeasonality.Range("M3").Value = 0.2
$wsSeasonality.Range("N3").Formula = "=SUM(B3:M3)"

# Row 4 - Jupiter.
$wsSeasonality.Range("B4").Value = 0
$wsSeasonality.Range("C4").Value = 0
$wsSeasonality.Range("D4").Value = 0
$wsSeasonality.Range("E4").Value = 0
$wsSeasonality.Range("F4").Value = 0
$wsSeasonality.Range("G4").Value = 0
$wsSeasonality.Range("H4").Value = 0
$wsSeasonality.Range("I4").Value = 0.1
$wsSeasonality.Range("J4").Value = 0.15
$wsSeasonality.Range("K4").Value = 0.25
$wsSeasonality.Range("L4").Value = 0.35
$wsSeasonality.Range("M4").Value = 0.15
$wsSeasonality.Range("N4").Formula = "=SUM(B4:M4)"

# Row 5 - Tesla (only Jan-Apr populated, rest left blank/formatted).
$wsSeasonality.Range("B5").Value = 0.35
$wsSeasonality.Range("C5").Value = 0.3
$wsSeasonality.Range("D5").Value = 0.25
$wsSeasonality.Range("E5").Value = 0.1
$wsSeasonality.Range("F5:M5").Value = ""
$wsSeasonality.Range("N5").Formula = "=SUM(B5:M5)"

# Row 6 - Venus: three literal "=1/6" formulas, then a shared-formula
# block (E6:G6), then H6:M6 left blank/formatted only.
$wsSeasonality.Range("B6").Formula = "=1/6"
$wsSeasonality.Range("C6").Formula = "=1/6"
$wsSeasonality.Range("D6").Formula = "=1/6"
$wsSeasonality.Range("E6:G6").Formula = "=1/6"
$wsSeasonality.Range("H6:M6").Value = ""
$wsSeasonality.Range("N6").Formula = "=SUM(B6:M6)"

# Rows 7-8: formatting-only cells, no values.
$wsSeasonality.Range("B7:M7").Value = ""
$wsSeasonality.Range("B8:M8").Value = ""

# Row 13 marker.
$wsSeasonality.Range("I13").Value = "x"

# ---------------------------------------------------------------------------
# 5. Selections per sheet + active sheet/tab (Project Fields stays active).
# ---------------------------------------------------------------------------
$wsAllocation.Range("B2").Select()
$wsSeasonality.Range("J13").Select()
$wsFields.Range("C9").Select()
$wsFields.Activate()
